$wb = $excel.ActiveWorkbook

# --- NameData (Sheet1): add new row 5 ---
$wsName = $wb.Worksheets.Item("NameData")
$wsName.Range("A5").Value = "4"
[void]$wsName.Range("A6").Select()

# --- AddressData (Sheet4): add new row 4 ---
$wsAddress = $wb.Worksheets.Item("AddressData")
$wsAddress.Range("A4").Value = "3"

# --- EmailAndPhoneData (Sheet5): add new rows 3-5, becomes the active tab ---
$wsEmail = $wb.Worksheets.Item("EmailAndPhoneData")
$wsEmail.Range("A3").Value = "2"
$wsEmail.Range("C3").Value = "2409876325"
$wsEmail.Range("A4").Value = "3"
$wsEmail.Range("A5").Value = "4"
$wsEmail.Range("B5").Value = "iahmed1@govolution.com"
$wsEmail.Range("B5").Style = "Normal"

# --- UDFData (Sheet6): add new row 3 ---
$wsUdf = $wb.Worksheets.Item("UDFData")
$wsUdf.Range("A3").Value = "2"
[void]$wsUdf.Range("B4").Select()

# --- AddressData selection now single cell A5 (after row 4 was filled) ---
[void]$wsAddress.Range("A5").Select()

# --- EmailAndPhoneData becomes the active sheet/tab, cursor at B5 ---
[void]$wsEmail.Activate()
[void]$wsEmail.Range("B5").Select()
